$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Modify Transaction" ---
# Insert a new "ReceiptNumber" row above the existing last row.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
$ws1.Rows.Item(4).Insert() | Out-Null
$ws1.Cells.Item(4,1).Value2 = "ReceiptNumber"
$ws1.Cells.Item(4,2).Value2 = 1213
$ws1.Range("C8").Select() | Out-Null

# --- Sheet 2: "Modify Transaction1" ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()
$ws2.Rows.Item(4).Insert() | Out-Null
$ws2.Cells.Item(4,1).Value2 = "ReceiptNumber"
$ws2.Cells.Item(4,2).Value2 = 123
$ws2.Range("C9").Select() | Out-Null

# --- Sheet 3: "Modify Transaction2" ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Activate()
$ws3.Rows.Item(4).Insert() | Out-Null
$ws3.Cells.Item(4,1).Value2 = "ReceiptNumber"
$ws3.Cells.Item(4,2).Value2 = 33
$ws3.Range("C9").Select() | Out-Null

# Restore the originally active sheet ("Transactions", sheet index 6 / activeTab=5)
$ws6 = $wb.Worksheets.Item(6)
$ws6.Activate()
